# Add daily power records for rows 86 and 87 (End Time values that were
# previously missing), plus a Start Time for row 87. This fills in the
# "Duration"/"Second Duration"/"Absolute Value" calculated columns via the
# table's existing shared formulas, and moves the active selection to the
# cell that was last edited (D87), matching Excel's usual post-edit focus.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 86: End Time (C86) was blank -> now recorded.
$ws.Range("C86").Value = 0.99930555555555556

# Row 87: Start Time (B87) and End Time (C87) were blank -> now recorded.
$ws.Range("B87").Value = 0
$ws.Range("C87").Value = 0.34513888888888888

# Move the selection to reflect the last cell touched by this edit.
$ws.Range("D87").Select()
